$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 287, shifting the existing rows 287-298
# down to 289-300 (dates/values/styles all move with the rows).
$ws.Rows("287:288").Insert()

# Row 287: new "Fukumoto" / "Provincia de Melipilla" record
$ws.Range("A287").Value = 11
$ws.Range("B287").Value = "Vega Monumental Concepción"
$ws.Range("C287").Value = "Bíobío"
$ws.Range("D287").Value = 44747
$ws.Range("E287").Value = 8
$ws.Range("F287").Value = "Fruta"
$ws.Range("G287").Value = 100102
$ws.Range("H287").Value = "Cítricos"
$ws.Range("I287").Value = 100102005
$ws.Range("J287").Value = "Naranja"
$ws.Range("K287").Value = "Fukumoto"
$ws.Range("L287").Value = "Primera"
$ws.Range("M287").Value = 100
$ws.Range("N287").Value = 7000
$ws.Range("O287").Value = 8000
$ws.Range("P287").Value = 7500
$ws.Range("Q287").Value = "$/caja 15 kilos empedrada"
$ws.Range("R287").Value = "Provincia de Melipilla"
$ws.Range("S287").Value = 500
$ws.Range("T287").Value = 15

# Row 288: new "Fukumoto" / "Provincia de Melipilla" record
$ws.Range("A288").Value = 11
$ws.Range("B288").Value = "Vega Monumental Concepción"
$ws.Range("C288").Value = "Bíobío"
$ws.Range("D288").Value = 44747
$ws.Range("E288").Value = 8
$ws.Range("F288").Value = "Fruta"
$ws.Range("G288").Value = 100102
$ws.Range("H288").Value = "Cítricos"
$ws.Range("I288").Value = 100102005
$ws.Range("J288").Value = "Naranja"
$ws.Range("K288").Value = "Fukumoto"
$ws.Range("L288").Value = "Segunda"
$ws.Range("M288").Value = 50
$ws.Range("N288").Value = 6000
$ws.Range("O288").Value = 6000
$ws.Range("P288").Value = 6000
$ws.Range("Q288").Value = "$/caja 15 kilos empedrada"
$ws.Range("R288").Value = "Provincia de Melipilla"
$ws.Range("S288").Value = 400
$ws.Range("T288").Value = 15
